$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value look like plain decimal numbers need to be pinned to
# Text format first, otherwise Excel auto-converts them to numeric cells and
# mangles the literal formatting (trailing zeros, leading zeros, etc.). The
# NumberFormat + ClearFormats sequence forces literal text while leaving the
# cell style back at the workbook default (matches the original inline-string cells).
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range("D2").Value = '26.178.60'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '1.600.47'
$ws.Range("E3").Value = '  -0.38%  '

$ws.Range("E4").Value = '  +0.41%  '

Set-TextValue "D5" '212.00'
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -0.95%  '

$ws.Range("E9").Value = '  -0.80%  '

Set-TextValue "D10" '18.17'
$ws.Range("E10").Value = '  -1.89%  '

Set-TextValue "D11" '0.0811'
$ws.Range("E11").Value = '  +2.58%  '

$ws.Range("D12").Value = '1.822.22'
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("D13").Value = '1.601.25'
$ws.Range("E13").Value = '  -0.32%  '

Set-TextValue "D14" '4.02'
$ws.Range("E14").Value = '  -1.11%  '

Set-TextValue "D15" '0.518'
$ws.Range("E15").Value = '  +0.65%  '

$ws.Range("D16").Value = '26.181.59'
$ws.Range("E16").Value = '  +0.08%  '

Set-TextValue "D17" '60.93'
$ws.Range("E17").Value = '  +0.10%  '

$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("E19").Value = '  +0.35%  '

Set-TextValue "D20" '204.21'
$ws.Range("E20").Value = '  +2.81%  '

$ws.Range("E21").Value = '  -0.11%  '

Set-TextValue "D22" '9.27'
$ws.Range("E22").Value = '  -2.44%  '

Set-TextValue "D23" '6.03'
$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("E24").Value = '  +12.22%  '

Set-TextValue "D25" '144.39'
$ws.Range("E25").Value = '  +1.56%  '

$ws.Range("E26").Value = '  +0.47%  '

$ws.Range("E27").Value = '  -7.68%  '

Set-TextValue "D28" '15.20'
$ws.Range("E28").Value = '  -0.56%  '

$ws.Range("E29").Value = '  +0.15%  '

Set-TextValue "D30" '0.0484'
$ws.Range("E30").Value = '  +1.45%  '

$ws.Range("E31").Value = '  -0.70%  '

$ws.Range("E32").Value = '  -0.32%  '

$ws.Range("E33").Value = '  -4.72%  '

$ws.Range("E34").Value = '  -2.12%  '

$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("D36").Value = '1.136.68'
$ws.Range("E36").Value = '  +2.58%  '

Set-TextValue "D37" '0.0163'
$ws.Range("E37").Value = '  +6.72%  '

$ws.Range("E38").Value = '  +0.42%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D39" '2.32'
$ws.Range("E39").Value = '  -1.64%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D40" '0.787'
$ws.Range("E40").Value = '  -0.80%  '

Set-TextValue "D41" '0.494'
$ws.Range("E41").Value = '  -2.93%  '

Set-TextValue "D42" '0.783'
$ws.Range("E42").Value = '  -2.35%  '

$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("D44").Value = '1.737.97'
$ws.Range("E44").Value = '  +0.16%  '

Set-TextValue "D45" '92.10'
$ws.Range("E45").Value = '  -0.86%  '

$ws.Range("E46").Value = '  -3.48%  '

Set-TextValue "D47" '54.09'
$ws.Range("E47").Value = '  +0.26%  '

$ws.Range("E48").Value = '  -0.31%  '

$ws.Range("E49").Value = '  -0.19%  '

$ws.Range("E50").Value = '  +0.61%  '

$ws.Range("D51").Value = '0.0₇0947'
$ws.Range("E51").Value = '  -11.64%  '

